$wb = $excel.ActiveWorkbook

# Navigate to the "Collections" worksheet and remove the two columns that
# held "Institutional Repository items" (D) and "Televison News Archive
# video segments" (E) data -- they shift "Linear feet in Special
# Collections & University Archives" left into column D.
$ws = $wb.Worksheets.Item("Collections")
$ws.Columns("D:E").Delete()

# The previously active sheet ("Visitors") loses its selection/tab focus
# once a different sheet is selected below.
$wsVisitors = $wb.Worksheets.Item("Visitors")
$wsVisitors.Range("B11").Select()

# Leave the selection where the editor last clicked, and make this sheet
# the active tab of the workbook.
$ws.Range("K17").Select()
$ws.Activate()
